# =====================================================================
# 260123 Insident and comment table create. make the ERD
#
# Adds two new DB-column-doc blocks to the bottom of the sheet:
#   - "insident"         (rows 46-56, PK column merged like the other tables)
#   - "insident_comment" (rows 57-62, not merged)
# and widens columns C/F so the new, longer column names/descriptions fit.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Stamp the new rows with the same look as the existing table blocks
#    before filling in any text, by copying formatting from row 3
#    (B = centered "PK column" style, C:F = plain bordered style).
# ---------------------------------------------------------------------
$ws.Range("B3:F3").Copy()
$ws.Range("B46:F56").PasteSpecial(-4122)

# Rows 57-62 ("insident_comment") are not merged in column B, so every
# column there (B:F) should use the plain bordered style from C3.
$ws.Range("C3:F3").Copy()
$ws.Range("C57:F62").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("B57:B62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Fill in the table/column documentation text
# ---------------------------------------------------------------------
$ws.Range("B46").Value2 = "insident"
$ws.Range("C46").Value2 = "insident_id"
$ws.Range("D46").Value2 = "v"
$ws.Range("E46").Value2 = "v"
$ws.Range("F46").Value2 = "Insident ID"

$ws.Range("C47").Value2 = "insident_title"
$ws.Range("E47").Value2 = "v"
$ws.Range("F47").Value2 = "Insident 제목"

$ws.Range("C48").Value2 = "insident_line_name"
$ws.Range("F48").Value2 = "인시던트 발생 호선명"

$ws.Range("C49").Value2 = "insident_station_id "
$ws.Range("E49").Value2 = "v"
$ws.Range("F49").Value2 = "인시던트 발생 역ID"

$ws.Range("C50").Value2 = "insident_station_name"
$ws.Range("F50").Value2 = "인시던트 발생 역 이름"

$ws.Range("C51").Value2 = "insident_content"
$ws.Range("F51").Value2 = "인시던트 발생 내용"

$ws.Range("C52").Value2 = "insident_status"
$ws.Range("E52").Value2 = "v"
$ws.Range("F52").Value2 = "인시던트 상태 (1:오픈 2:해결중 3:완료 4:비활성화)"

$ws.Range("C53").Value2 = "user_id"
$ws.Range("E53").Value2 = "v"
$ws.Range("F53").Value2 = "작성자 ID ( FK storage.user(user_id) )"

$ws.Range("C54").Value2 = "user_name"
$ws.Range("F54").Value2 = "작성자 이름"

$ws.Range("C55").Value2 = "create_at"
$ws.Range("E55").Value2 = "v"
$ws.Range("F55").Value2 = "등록일자 Default Current_Timestamp"

$ws.Range("C56").Value2 = "complete_at"
$ws.Range("F56").Value2 = "완료일자"

$ws.Range("B57").Value2 = "insident_comment"
$ws.Range("C57").Value2 = "comment_id"
$ws.Range("D57").Value2 = "v"
$ws.Range("E57").Value2 = "v"
$ws.Range("F57").Value2 = "댓글 ID"

$ws.Range("C58").Value2 = "insident_id"
$ws.Range("E58").Value2 = "v"
$ws.Range("F58").Value2 = "인시던트 ID ( FK storage.insident(insident_id) )"

$ws.Range("C59").Value2 = "comment_content"
$ws.Range("E59").Value2 = "v"
$ws.Range("F59").Value2 = "댓글 내용"

$ws.Range("C60").Value2 = "user_id"
$ws.Range("E60").Value2 = "v"
$ws.Range("F60").Value2 = "작성자 ID ( FK storage.user(user_id) )"

$ws.Range("C61").Value2 = "user_name"
$ws.Range("F61").Value2 = "작성자 이름"

$ws.Range("C62").Value2 = "create_at"
$ws.Range("E62").Value2 = "v"
$ws.Range("F62").Value2 = "등록일자 Default Current_Timestamp"

# ---------------------------------------------------------------------
# 3) Merge the "insident" table PK column (B46:B56), matching the other
#    table blocks above it, then drop the inside horizontal edges so the
#    merged block only shows an outer box (thin top on row46, thin
#    bottom on row56, none on the rows between).
# ---------------------------------------------------------------------
$ws.Range("B46:B56").Merge()
$ws.Range("B46:B56").Borders.Item(12).LineStyle = -4142

# ---------------------------------------------------------------------
# 4) Column C/F grew to fit the new, longer column names & descriptions
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 20.684151785714285
$ws.Columns.Item(6).ColumnWidth = 46.184151785714285

# ---------------------------------------------------------------------
# 5) Leave the selection where the author left it: scrolled down onto
#    the newly-added "insident_comment" block.
# ---------------------------------------------------------------------
$ws.Range("I60").Select()
